# This script permutes the per-trial metadata columns (H,I,K,L,M,N,O,P,Q,R,S,T,U,V)
# across data rows 2-41 of the sheet. Row 32 is left untouched. Columns A-G and J
# (subject_id, task, block_total, block_scene, trial_block, trial_total, target_cat,
# cond_mem) are not touched either; only the category/condition/answer/stimulus/
# rating columns get reassigned to a different row's values, per a fixed permutation
# (20 unique input-file variants get generated from the same base data by shuffling
# the trial rows; here we only store the resulting cell values).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column numbers for H, I, K, L, M, N, O, P, Q, R, S, T, U, V
$cols = @(8, 9, 11, 12, 13, 14, 15, 16, 17, 18, 19, 20, 21, 22)

# Mapping of destination row -> source row (values are copied FROM source TO destination).
# Row 32 is intentionally absent (unchanged).
$mapping = @{
    2  = 17
    3  = 16
    4  = 12
    5  = 40
    6  = 39
    7  = 24
    8  = 5
    9  = 20
    10 = 15
    11 = 3
    12 = 4
    13 = 31
    14 = 36
    15 = 29
    16 = 33
    17 = 34
    18 = 22
    19 = 18
    20 = 10
    21 = 13
    22 = 30
    23 = 25
    24 = 35
    25 = 41
    26 = 38
    27 = 23
    28 = 21
    29 = 2
    30 = 11
    31 = 14
    33 = 8
    34 = 26
    35 = 19
    36 = 28
    37 = 7
    38 = 27
    39 = 6
    40 = 9
    41 = 37
}

# Snapshot the current (pre-edit) values for every row 2-41 and every tracked column,
# so the writes below never read already-overwritten data.
$snapshot = @{}
for ($r = 2; $r -le 41; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value()
    }
    $snapshot[$r] = $rowVals
}

# Apply the permutation: destination row gets the source row's snapshot values.
foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    $srcVals = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $ws.Cells.Item($destRow, $c).Value = $srcVals[$c]
    }
}
